# Update gh-pages to output generated at 456a3b4
# Applies numeric "想去人数" (interest count) refreshes across the four
# sheets, one ticket-status flip, and a text re-sync of two rows in the
# "全部类型" (all-types) aggregation sheet.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item(1)   # 展览
$sheetShow    = $wb.Worksheets.Item(2)   # 演出
$sheetLocal   = $wb.Worksheets.Item(3)   # 本地生活
$sheetAll     = $wb.Worksheets.Item(4)   # 全部类型

# ---------------------------------------------------------------------
# 展览 (sheet 1) — column F ("想去人数") updates
# ---------------------------------------------------------------------
$sheetExhibit.Range("F2").Value  = 1925
$sheetExhibit.Range("F5").Value  = 401
$sheetExhibit.Range("F7").Value  = 847
$sheetExhibit.Range("F8").Value  = 1238
$sheetExhibit.Range("F10").Value = 346
$sheetExhibit.Range("F11").Value = 125
$sheetExhibit.Range("F12").Value = 2482
$sheetExhibit.Range("F14").Value = 340
$sheetExhibit.Range("F18").Value = 45
$sheetExhibit.Range("F19").Value = 1505
$sheetExhibit.Range("F20").Value = 404953
$sheetExhibit.Range("F21").Value = 1203
$sheetExhibit.Range("F22").Value = 151
$sheetExhibit.Range("F23").Value = 0
$sheetExhibit.Range("F24").Value = 1365
$sheetExhibit.Range("F25").Value = 1349
$sheetExhibit.Range("F26").Value = 941
$sheetExhibit.Range("F28").Value = 1282
$sheetExhibit.Range("F29").Value = 170
$sheetExhibit.Range("F30").Value = 1238
$sheetExhibit.Range("F31").Value = 408
$sheetExhibit.Range("F32").Value = 128
$sheetExhibit.Range("F33").Value = 939
$sheetExhibit.Range("F35").Value = 1800
$sheetExhibit.Range("F36").Value = 433
$sheetExhibit.Range("F37").Value = 29
$sheetExhibit.Range("F38").Value = 143
$sheetExhibit.Range("F39").Value = 14
$sheetExhibit.Range("F40").Value = 2210
$sheetExhibit.Range("F41").Value = 122
$sheetExhibit.Range("F43").Value = 2570

# ---------------------------------------------------------------------
# 演出 (sheet 2) — column F updates + one G (ticket price -> sold out)
# ---------------------------------------------------------------------
$sheetShow.Range("F5").Value  = 54
$sheetShow.Range("G6").Value  = "不可售"
$sheetShow.Range("F13").Value = 0
$sheetShow.Range("F17").Value = 54
$sheetShow.Range("F18").Value = 54
$sheetShow.Range("F20").Value = 277
$sheetShow.Range("F22").Value = 263
$sheetShow.Range("F26").Value = 50
$sheetShow.Range("F27").Value = 50
$sheetShow.Range("F29").Value = 40
$sheetShow.Range("F30").Value = 208
$sheetShow.Range("F32").Value = 41
$sheetShow.Range("F33").Value = 16
$sheetShow.Range("F34").Value = 81
$sheetShow.Range("F36").Value = 156
$sheetShow.Range("F37").Value = 74

# ---------------------------------------------------------------------
# 本地生活 (sheet 3) — column F updates
# ---------------------------------------------------------------------
$sheetLocal.Range("F5").Value  = 2994
$sheetLocal.Range("F6").Value  = 4817
$sheetLocal.Range("F9").Value  = 643
$sheetLocal.Range("F10").Value = 891
$sheetLocal.Range("F11").Value = 522
$sheetLocal.Range("F12").Value = 566
$sheetLocal.Range("F13").Value = 1295
$sheetLocal.Range("F14").Value = 367
$sheetLocal.Range("F15").Value = 1065

# ---------------------------------------------------------------------
# 全部类型 (sheet 4) — column F updates
# ---------------------------------------------------------------------
$sheetAll.Range("F2").Value  = 1925
$sheetAll.Range("F5").Value  = 4817
$sheetAll.Range("F7").Value  = 643
$sheetAll.Range("F8").Value  = 891
$sheetAll.Range("F9").Value  = 522
$sheetAll.Range("F11").Value = 566
$sheetAll.Range("F12").Value = 1295
$sheetAll.Range("F13").Value = 401

# Rows 14 & 15 on this sheet were stale copies of two different 06-29
# carnival events; resync every field (C..I) to match the refreshed
# listing (same events now shown on the 展览 sheet as rows 6 & 7).
$sheetAll.Range("C14").Value = "上海·创世次元动漫游戏嘉年华3.0"
$sheetAll.Range("D14").Value = "中环立交桥苏宁天御国际广场西南侧约240米 轮客行轮滑馆(普陀店)"
$sheetAll.Range("E14").Value = "2024.06.29 10:00-06.30 17:00"
$sheetAll.Range("F14").Value = 1788
$sheetAll.Range("G14").Value = 58
$sheetAll.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=86506"
$sheetAll.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202405/Clkfdwic1716894666596.jpeg"

$sheetAll.Range("C15").Value = "上海·第六十三届燃梦星辰动漫嘉年华"
$sheetAll.Range("D15").Value = "陆宝山路155号 佘山·旭辉里"
$sheetAll.Range("E15").Value = "2024.06.29 14:00-06.29 18:00"
$sheetAll.Range("F15").Value = 847
$sheetAll.Range("G15").Value = 58.8
$sheetAll.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=85231"
$sheetAll.Range("I15").Value = "//i2.hdslb.com/bfs/openplatform/202406/Tx1g80PC1717667546414.jpeg"

$sheetAll.Range("F16").Value = 1238
$sheetAll.Range("F17").Value = 346
$sheetAll.Range("F18").Value = 1065
$sheetAll.Range("F19").Value = 2482
$sheetAll.Range("F21").Value = 340
$sheetAll.Range("F25").Value = 1505
$sheetAll.Range("F28").Value = 1203
$sheetAll.Range("F29").Value = 151
$sheetAll.Range("F31").Value = 1349
$sheetAll.Range("F32").Value = 941
$sheetAll.Range("F33").Value = 1282
$sheetAll.Range("F34").Value = 170
$sheetAll.Range("F36").Value = 54
$sheetAll.Range("F37").Value = 1238
$sheetAll.Range("F38").Value = 408
$sheetAll.Range("F39").Value = 939
$sheetAll.Range("F41").Value = 1800
$sheetAll.Range("F42").Value = 50
$sheetAll.Range("F43").Value = 29
$sheetAll.Range("F44").Value = 143
$sheetAll.Range("F45").Value = 2210
$sheetAll.Range("F46").Value = 122
$sheetAll.Range("F48").Value = 2570
$sheetAll.Range("F51").Value = 74
